$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.390674242509391
$ws.Range("C2").Value = 0.05422536076038398
$ws.Range("D2").Value = 0.03938223372932015
$ws.Range("F2").Value = 0.7997251448389093
$ws.Range("G2").Value = 0.6436838506097615
$ws.Range("H2").Value = 0.7644538316676517
$ws.Range("K2").Value = 0.3682658822932865
$ws.Range("B3").Value = 0.34940106566998
$ws.Range("C3").Value = 0.04868489409926724
$ws.Range("D3").Value = 0.03721867884566166
$ws.Range("F3").Value = 0.8040265888460212
$ws.Range("G3").Value = 0.6498349275094171
$ws.Range("H3").Value = 0.7723479790848273
$ws.Range("K3").Value = 0.3241280516861309
$ws.Range("B4").Value = 0.3241034189544791
$ws.Range("C4").Value = 0.04526212056113366
$ws.Range("D4").Value = 0.03588131487120449
$ws.Range("F4").Value = 0.8072886745315131
$ws.Range("G4").Value = 0.6541833909715749
$ws.Range("H4").Value = 0.777627473804138
$ws.Range("K4").Value = 0.2970092953288486
$ws.Range("B5").Value = 0.3138059892730496
$ws.Range("C5").Value = 0.04386211876133927
$ws.Range("D5").Value = 0.03533412012694725
$ws.Range("F5").Value = 0.8087739092059678
$ws.Range("G5").Value = 0.6560988264466303
$ws.Range("H5").Value = 0.7798875928337381
$ws.Range("K5").Value = 0.2859541299455941
$ws.Range("B6").Value = 0.3120968205477936
$ws.Range("C6").Value = 0.04362933770096333
$ws.Range("D6").Value = 0.03524312662468532
$ws.Range("F6").Value = 0.8090299406456651
$ws.Range("G6").Value = 0.6564255338128646
$ws.Range("H6").Value = 0.7802694461752608
$ws.Range("K6").Value = 0.2841181971145943
$ws.Range("B7").Value = 0.3239644967979984
$ws.Range("C7").Value = 0.04524326059227235
$ws.Range("D7").Value = 0.03587394410090639
$ws.Range("F7").Value = 0.8073080739827887
$ws.Range("G7").Value = 0.6542086430456564
$ws.Range("H7").Value = 0.7776575145952265
$ws.Range("K7").Value = 0.2968602172929309
$ws.Range("B8").Value = 0.3764343167666482
$ws.Range("C8").Value = 0.05231938251074553
$ws.Range("D8").Value = 0.03863811604720979
$ws.Range("F8").Value = 0.8010792509258238
$ws.Range("G8").Value = 0.6456858786288535
$ws.Range("H8").Value = 0.7670859399584771
$ws.Range("K8").Value = 0.3530511462171262
$ws.Range("B9").Value = 0.4796645823945482
$ws.Range("C9").Value = 0.0660276655693508
$ws.Range("D9").Value = 0.04398636047582727
$ws.Range("F9").Value = 0.7938034284486974
$ws.Range("G9").Value = 0.6335249151098594
$ws.Range("H9").Value = 0.7497893180849786
$ws.Range("K9").Value = 0.4630842188480528
$ws.Range("B10").Value = 0.5557023852577174
$ws.Range("C10").Value = 0.07599488843897007
$ws.Range("D10").Value = 0.04787013156047237
$ws.Range("F10").Value = 0.7914859573966169
$ws.Range("G10").Value = 0.6273881353054094
$ws.Range("H10").Value = 0.7391790820247763
$ws.Range("K10").Value = 0.5438180604407421
$ws.Range("B11").Value = 0.5903343865468287
$ws.Range("C11").Value = 0.08050627471426708
$ws.Range("D11").Value = 0.04962676991370074
$ws.Range("F11").Value = 0.7910929495405767
$ws.Range("G11").Value = 0.6252087969555902
$ws.Range("H11").Value = 0.7348085641337576
$ws.Range("K11").Value = 0.5805207269119421
$ws.Range("B12").Value = 0.6034543429043424
$ws.Range("C12").Value = 0.08221129930640814
$ws.Range("D12").Value = 0.05029047748327997
$ws.Range("F12").Value = 0.7910394943906525
$ws.Range("G12").Value = 0.6244719787945172
$ws.Range("H12").Value = 0.7332192314507324
$ws.Range("K12").Value = 0.5944153426104322
$ws.Range("B13").Value = 0.600628485382714
$ws.Range("C13").Value = 0.08184424134724111
$ws.Range("D13").Value = 0.0501476032218946
$ws.Range("F13").Value = 0.7910467608491203
$ws.Range("G13").Value = 0.6246267255791764
$ws.Range("H13").Value = 0.7335585996721932
$ws.Range("K13").Value = 0.5914230693223601
$ws.Range("B14").Value = 0.5914136643432641
$ws.Range("C14").Value = 0.0806466154268719
$ws.Range("D14").Value = 0.04968140368735163
$ws.Range("F14").Value = 0.7910866387989017
$ws.Range("G14").Value = 0.6251464035273244
$ws.Range("H14").Value = 0.7346764917482602
$ws.Range("K14").Value = 0.5816639260661418
$ws.Range("B15").Value = 0.5857700316316539
$ws.Range("C15").Value = 0.07991259807727147
$ws.Range("D15").Value = 0.04939564737705382
$ws.Range("F15").Value = 0.7911234934212601
$ws.Range("G15").Value = 0.6254762519184851
$ws.Range("H15").Value = 0.7353697902471055
$ws.Range("K15").Value = 0.5756856480053898
$ws.Range("B16").Value = 0.5534399075717999
$ws.Range("C16").Value = 0.07569959508747104
$ws.Range("D16").Value = 0.04775512430211393
$ws.Range("F16").Value = 0.7915249712984433
$ws.Range("G16").Value = 0.6275429157379193
$ws.Range("H16").Value = 0.7394738924248543
$ws.Range("K16").Value = 0.5414189379648633
$ws.Range("B17").Value = 0.5336168487352495
$ws.Range("C17").Value = 0.07310917822303509
$ws.Range("D17").Value = 0.04674609912182603
$ws.Range("F17").Value = 0.7919408233522702
$ws.Range("G17").Value = 0.6289678603771023
$ws.Range("H17").Value = 0.7421085210468874
$ws.Range("K17").Value = 0.5203910298473602
$ws.Range("B18").Value = 0.5222191331659189
$ws.Range("C18").Value = 0.07161710129318521
$ws.Range("D18").Value = 0.04616478539729485
$ws.Range("F18").Value = 0.7922422307694958
$ws.Range("G18").Value = 0.629845060878381
$ws.Range("H18").Value = 0.7436668172524321
$ws.Range("K18").Value = 0.5082941450152987
$ws.Range("B19").Value = 0.5183607611832031
$ws.Range("C19").Value = 0.07111154460844205
$ws.Range("D19").Value = 0.04596780076087725
$ws.Range("F19").Value = 0.7923549598145243
$ws.Range("G19").Value = 0.6301519482616555
$ws.Range("H19").Value = 0.7442017994781338
$ws.Range("K19").Value = 0.5041979870999853
$ws.Range("B20").Value = 0.5357266383466879
$ws.Range("C20").Value = 0.0733851542526196
$ws.Range("D20").Value = 0.04685361009511269
$ws.Range("F20").Value = 0.7918901135466925
$ws.Range("G20").Value = 0.62881020726158
$ws.Range("H20").Value = 0.7418236167960117
$ws.Range("K20").Value = 0.5226297166686038
$ws.Range("B21").Value = 0.5941201332717299
$ws.Range("C21").Value = 0.08099847835771357
$ws.Range("D21").Value = 0.04981837863095251
$ws.Range("F21").Value = 0.7910723351400009
$ws.Range("G21").Value = 0.6249913578295008
$ws.Range("H21").Value = 0.7343463562366281
$ws.Range("K21").Value = 0.5845305324730532
$ws.Range("B22").Value = 0.6323159009121468
$ws.Range("C22").Value = 0.08595472985801678
$ws.Range("D22").Value = 0.05174730437641983
$ws.Range("F22").Value = 0.7910938810094947
$ws.Range("G22").Value = 0.6230112576416218
$ws.Range("H22").Value = 0.7298424565133956
$ws.Range("K22").Value = 0.6249634052894919
$ws.Range("B23").Value = 0.6119273064155948
$ws.Range("C23").Value = 0.0833112911777647
$ws.Range("D23").Value = 0.05071861143703416
$ws.Range("F23").Value = 0.7910314125810345
$ws.Range("G23").Value = 0.6240207524871124
$ws.Range("H23").Value = 0.7322112032885855
$ws.Range("K23").Value = 0.6033858784408892
$ws.Range("B24").Value = 0.5347728060896202
$ws.Range("C24").Value = 0.07326039424326325
$ws.Range("D24").Value = 0.04680500815535993
$ws.Range("F24").Value = 0.791912845320816
$ws.Range("G24").Value = 0.6288813016571737
$ws.Range("H24").Value = 0.7419522861605827
$ws.Range("K24").Value = 0.5216176301826465
$ws.Range("B25").Value = 0.4517031184913094
$ws.Range("C25").Value = 0.06233738850920645
$ws.Range("D25").Value = 0.04254741605542023
$ws.Range("F25").Value = 0.7952412062761027
$ws.Range("G25").Value = 0.6363250802329929
$ws.Range("H25").Value = 0.7541004433392402
$ws.Range("K25").Value = 0.4333354016445412
